$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Write the new contact names in the same order the unique strings were first
# introduced, so the shared-strings table is built with the expected ordering.
$ws.Range("A3").Value = "Luis"
$ws.Range("A5").Value = "Maria Helena"
$ws.Range("A4").Value = "Bruno"
$ws.Range("A6").Value = "Bruninho"
$ws.Range("A7").Value = "João"
$ws.Range("A8").Value = "Aquino"
$ws.Range("A2").Value = "Marcelo"

# Move the active selection to A2, mirroring the updated sheet view
$ws.Range("A2").Select()
